$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.806.02'
$ws.Range("E2").Value = '  -0.49%  '

$ws.Range("D3").Value = '2.497.37'
$ws.Range("E3").Value = '  +2.47%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '535.27'
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.55'
$ws.Range("E6").Value = '  -3.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("D9").Value = '2.524.12'
$ws.Range("E9").Value = '  +3.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0992'
$ws.Range("E10").Value = '  +0.86%  '

$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.46'
$ws.Range("E12").Value = '  +2.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.349'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").Value = '2.936.90'
$ws.Range("E14").Value = '  +2.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.36'
$ws.Range("E15").Value = '  -2.62%  '

$ws.Range("D16").Value = '58.741.73'
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("E17").Value = '  +1.04%  '

$ws.Range("D18").Value = '2.511.25'
$ws.Range("E18").Value = '  +1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.16'
$ws.Range("E19").Value = '  +1.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.99'
$ws.Range("E21").Value = '  -0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +3.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.73'
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.66'
$ws.Range("E24").Value = '  +2.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.434'
$ws.Range("E25").Value = '  -6.81%  '

$ws.Range("E26").Value = '  +1.38%  '

$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.610.43'
$ws.Range("E27").Value = '  +1.95%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.994'
$ws.Range("E28").Value = '  +1.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.67'
$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.64'
$ws.Range("E30").Value = '  -2.34%  '

$ws.Range("D31").Value = '0.0₃0763'
$ws.Range("E31").Value = '  -0.84%  '

$ws.Range("E32").Value = '  -1.37%  '

$ws.Range("E33").Value = '  -8.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  -0.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.63'
$ws.Range("E35").Value = '  +0.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.42'
$ws.Range("E36").Value = '  +6.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.52'
$ws.Range("E37").Value = '  +1.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.30'
$ws.Range("E38").Value = '  -4.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.60'
$ws.Range("E39").Value = '  -6.55%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.54'
$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("E41").Value = '  -4.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '296.87'
$ws.Range("E42").Value = '  -4.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.62'
$ws.Range("E43").Value = '  -2.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.805'
$ws.Range("E44").Value = '  -4.91%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.995'
$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.602'
$ws.Range("E46").Value = '  +4.39%  '

$ws.Range("E47").Value = '  +0.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.59'
$ws.Range("E48").Value = '  +4.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0925'
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.50'
$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0226'
$ws.Range("E51").Value = '  -0.45%  '
